$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(5).Insert()

$ws.Range("E1").Value = "City"
$ws.Range("E2").Value = "København"
$ws.Range("E3").Value = "Malmö"
$ws.Range("E4").Value = "Oslo"

$ws.Range("E4").Select() | Out-Null
